$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.887.00"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.87"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.67"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5031"
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2569"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06397"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.74"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.659.21"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.277"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.863.92"
$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5465"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7926"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.51"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.884.89"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.92"
$ws.Range("E20").Value = "  -3.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.387"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.933"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.986"
$ws.Range("E23").Value = "  -1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.916"
$ws.Range("E25").Value = "  +9.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.00"

$ws.Range("E27").Value = "  -3.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.71"
$ws.Range("E28").Value = "  -0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.734"
$ws.Range("E29").Value = "  -3.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.245"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04944"
$ws.Range("E31").Value = "  -3.73%  "

$ws.Range("E32").Value = "  -2.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.198"
$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.544"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("E35").Value = "  +0.71%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.631"
$ws.Range("E36").Value = "  -4.16%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8938"
$ws.Range("E37").Value = "  -3.35%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.157.36"
$ws.Range("E38").Value = "  -0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5602"
$ws.Range("E39").Value = "  -2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01564"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.699"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8072"
$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.80"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.775.72"
$ws.Range("E45").Value = "  -0.79%  "

$ws.Range("E46").Value = "  -0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4528"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.04"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05064"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  -0.17%  "
